$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.504.19'
$ws.Range('E2').Value = '  -1.53%  '
$ws.Range('D3').Value = '2.163.18'
$ws.Range('E3').Value = '  -2.73%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '237.58'
$ws.Range('E5').Value = '  -2.28%  '
$ws.Range('D6').Value = '0.608'
$ws.Range('E6').Value = '  -3.05%  '
$ws.Range('D7').Value = '71.78'
$ws.Range('E7').Value = '  -3.26%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').Value = '0.571'
$ws.Range('E9').Value = '  -5.32%  '
$ws.Range('D10').Value = '39.58'
$ws.Range('E10').Value = '  -7.19%  '
$ws.Range('D11').Value = '0.0904'
$ws.Range('E11').Value = '  -5.93%  '
$ws.Range('D12').Value = '54.40'
$ws.Range('E12').Value = '  -3.90%  '
$ws.Range('E13').Value = '  -3.27%  '
$ws.Range('D14').Value = '6.64'
$ws.Range('E14').Value = '  -4.75%  '
$ws.Range('D15').Value = '2.487.67'
$ws.Range('E15').Value = '  -2.67%  '
$ws.Range('E16').Value = '  -0.27%  '
$ws.Range('D17').Value = '2.163.93'
$ws.Range('E17').Value = '  -2.48%  '
$ws.Range('D18').Value = '0.774'
$ws.Range('E18').Value = '  -7.39%  '
$ws.Range('D19').Value = '41.422.11'
$ws.Range('E19').Value = '  -1.27%  '
$ws.Range('E20').Value = '  -2.86%  '
$ws.Range('D21').Value = '69.93'
$ws.Range('E21').Value = '  -4.03%  '
$ws.Range('D22').Value = '5.77'
$ws.Range('E22').Value = '  -7.18%  '
$ws.Range('D23').Value = '9.93'
$ws.Range('E23').Value = '  -10.27%  '
$ws.Range('D24').Value = '226.03'
$ws.Range('E24').Value = '  -1.93%  '
$ws.Range('D25').Value = '2.01'
$ws.Range('E25').Value = '  -3.82%  '
$ws.Range('E26').Value = '  -0.15%  '
$ws.Range('D27').Value = '10.70'
$ws.Range('E28').Value = '  -9.74%  '
$ws.Range('D29').Value = '2.19'
$ws.Range('E29').Value = '  -3.67%  '
$ws.Range('D30').Value = '2.16'
$ws.Range('E30').Value = '  -1.82%  '
$ws.Range('D31').Value = '170.72'
$ws.Range('E31').Value = '  +2.50%  '
$ws.Range('D32').Value = '19.77'
$ws.Range('E32').Value = '  -3.97%  '
$ws.Range('D33').Value = '32.88'
$ws.Range('E33').Value = '  +9.69%  '
$ws.Range('E34').Value = '  -4.13%  '
$ws.Range('D35').Value = '5.33'
$ws.Range('E35').Value = '  -5.56%  '
$ws.Range('E36').Value = '  -3.68%  '
$ws.Range('D37').Value = '4.27'
$ws.Range('E37').Value = '  -1.11%  '
$ws.Range('E38').Value = '  -7.38%  '
$ws.Range('E39').Value = '  -0.23%  '
$ws.Range('D40').Value = '11.97'
$ws.Range('E40').Value = '  -9.47%  '
$ws.Range('E41').Value = '  -1.77%  '
$ws.Range('D42').Value = '5.34'
$ws.Range('E42').Value = '  -6.22%  '
$ws.Range('D43').Value = '58.89'
$ws.Range('E43').Value = '  -9.39%  '
$ws.Range('E44').Value = '  -2.79%  '
$ws.Range('E45').Value = '  -5.75%  '
$ws.Range('E46').Value = '  -3.81%  '
$ws.Range('D47').Value = '96.82'
$ws.Range('E47').Value = '  -7.16%  '
$ws.Range('D48').Value = '1.07'
$ws.Range('E48').Value = '  -3.80%  '
$ws.Range('E49').Value = '  -4.96%  '
$ws.Range('E50').Value = '  -7.58%  '
$ws.Range('E51').Value = '  -2.45%  '
